# Updated cryptos list on Fri Nov  3 23:35:46 UTC 2023 with GitHub Actions
#
# The sheet is a scraped crypto price table; every data row packs Coin (B),
# Link (C), Price (D) and 1h Volume change (E) as plain text cells so the
# original formatting (trailing zeros like "40.00", dotted thousands like
# "34.927.04", padded "  +1.20%  ") survives untouched. This refresh just
# rewrites the handful of cells the latest scrape changed - mostly D/E
# price+change pairs, plus a few rows (14-16) whose coin ordering shifted.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "34.927.04" },
    @{ Cell = "E2"; Value = "  -0.40%  " },
    @{ Cell = "D3"; Value = "1.843.27" },
    @{ Cell = "E3"; Value = "  +1.90%  " },
    @{ Cell = "E4"; Value = "  +0.00%  " },
    @{ Cell = "D5"; Value = "231.78" },
    @{ Cell = "E5"; Value = "  -0.30%  " },
    @{ Cell = "D6"; Value = "0.619" },
    @{ Cell = "E6"; Value = "  +1.20%  " },
    @{ Cell = "E7"; Value = "  +0.00%  " },
    @{ Cell = "D8"; Value = "40.00" },
    @{ Cell = "E8"; Value = "  -0.94%  " },
    @{ Cell = "D9"; Value = "0.330" },
    @{ Cell = "E9"; Value = "  +0.99%  " },
    @{ Cell = "D10"; Value = "0.0687" },
    @{ Cell = "E10"; Value = "  +0.39%  " },
    @{ Cell = "D11"; Value = "0.0981" },
    @{ Cell = "E11"; Value = "  -1.89%  " },
    @{ Cell = "D12"; Value = "2.109.54" },
    @{ Cell = "E12"; Value = "  +1.86%  " },
    @{ Cell = "D13"; Value = "11.61" },
    @{ Cell = "E13"; Value = "  +5.00%  " },
    @{ Cell = "B14"; Value = "Polygon" },
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic" },
    @{ Cell = "D14"; Value = "0.675" },
    @{ Cell = "E14"; Value = "  +1.80%  " },
    @{ Cell = "B15"; Value = "Polkadot" },
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot" },
    @{ Cell = "D15"; Value = "4.65" },
    @{ Cell = "E15"; Value = "  -0.24%  " },
    @{ Cell = "B16"; Value = "WrappedEther" },
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" },
    @{ Cell = "D16"; Value = "1.783.27" },
    @{ Cell = "E16"; Value = "  -1.54%  " },
    @{ Cell = "D17"; Value = "34.899.89" },
    @{ Cell = "E17"; Value = "  -0.34%  " },
    @{ Cell = "D18"; Value = "69.90" },
    @{ Cell = "E18"; Value = "  +0.21%  " },
    @{ Cell = "D19"; Value = "0.0₃0788" },
    @{ Cell = "E19"; Value = "  -0.24%  " },
    @{ Cell = "D20"; Value = "240.14" },
    @{ Cell = "E20"; Value = "  +1.02%  " },
    @{ Cell = "D21"; Value = "12.18" },
    @{ Cell = "E21"; Value = "  +2.03%  " },
    @{ Cell = "D22"; Value = "4.69" },
    @{ Cell = "E22"; Value = "  -0.79%  " },
    @{ Cell = "E23"; Value = "  +0.08%  " },
    @{ Cell = "E24"; Value = "  +0.97%  " },
    @{ Cell = "D25"; Value = "171.07" },
    @{ Cell = "E25"; Value = "  -0.51%  " },
    @{ Cell = "D26"; Value = "7.82" },
    @{ Cell = "E26"; Value = "  -0.95%  " },
    @{ Cell = "D27"; Value = "17.49" },
    @{ Cell = "E27"; Value = "  -0.13%  " },
    @{ Cell = "E28"; Value = "  +2.10%  " },
    @{ Cell = "E29"; Value = "  -3.03%  " },
    @{ Cell = "E30"; Value = "  +0.09%  " },
    @{ Cell = "D31"; Value = "0.0552" },
    @{ Cell = "E31"; Value = "  -1.06%  " },
    @{ Cell = "D32"; Value = "3.95" },
    @{ Cell = "E32"; Value = "  -4.49%  " },
    @{ Cell = "D33"; Value = "3.97" },
    @{ Cell = "E33"; Value = "  -1.41%  " },
    @{ Cell = "D34"; Value = "1.92" },
    @{ Cell = "E34"; Value = "  +8.56%  " },
    @{ Cell = "E35"; Value = "  +6.78%  " },
    @{ Cell = "D36"; Value = "1.47" },
    @{ Cell = "E36"; Value = "  +14.35%  " },
    @{ Cell = "D37"; Value = "0.699" },
    @{ Cell = "E37"; Value = "  +0.37%  " },
    @{ Cell = "D38"; Value = "1.08" },
    @{ Cell = "E38"; Value = "  +8.36%  " },
    @{ Cell = "D39"; Value = "90.51" },
    @{ Cell = "E39"; Value = "  -1.70%  " },
    @{ Cell = "D40"; Value = "1.350.77" },
    @{ Cell = "E40"; Value = "  +2.83%  " },
    @{ Cell = "E41"; Value = "  +0.39%  " },
    @{ Cell = "D42"; Value = "14.82" },
    @{ Cell = "E42"; Value = "  +2.49%  " },
    @{ Cell = "E43"; Value = "  +2.07%  " },
    @{ Cell = "E44"; Value = "  -2.58%  " },
    @{ Cell = "E45"; Value = "  +0.08%  " },
    @{ Cell = "E46"; Value = "  +2.33%  " },
    @{ Cell = "E47"; Value = "  -0.84%  " },
    @{ Cell = "D48"; Value = "2.023.05" },
    @{ Cell = "E48"; Value = "  +1.75%  " },
    @{ Cell = "D49"; Value = "3.46" },
    @{ Cell = "E49"; Value = "  +22.58%  " },
    @{ Cell = "D50"; Value = "1.01" },
    @{ Cell = "E50"; Value = "  +0.12%  " },
    @{ Cell = "D51"; Value = "0.0667" },
    @{ Cell = "E51"; Value = "  +0.48%  " }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)

    if ($u.Value -match '^[0-9]+(\.[0-9]+)?$') {
        # Looks like a plain number (e.g. "40.00", "11.61") - force the
        # cell to Text first so Excel keeps the literal digits/trailing
        # zeros instead of silently coercing to a numeric value.
        $range.NumberFormat = "@"
        $range.Value = $u.Value
        $range.Style = "Normal"
    } else {
        # Coin names, URLs, multi-dot prices ("34.927.04") and padded
        # percent strings ("  +1.20%  ") are never number-like, so a plain
        # assignment is safe and keeps cell styling untouched.
        $range.Value = $u.Value
    }
}
